$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-epoch validation accuracy values (column B) for the re-run
# with Token Embeddings + Decoder layer 12 frozen (M14).
$newValues = @{
    "2" = 0.453125
    "3" = 0.40625
    "4" = 0.375
    "6" = 0.296875
    "7" = 0.28125
    "8" = 0.328125
    "9" = 0.3125
    "10" = 0.328125
    "11" = 0.28125
    "12" = 0.296875
    "14" = 0.3125
    "15" = 0.328125
    "16" = 0.265625
    "17" = 0.296875
    "18" = 0.265625
    "19" = 0.265625
    "20" = 0.265625
    "21" = 0.296875
    "22" = 0.296875
    "23" = 0.25
    "25" = 0.234375
    "26" = 0.234375
    "28" = 0.21875
    "29" = 0.21875
    "30" = 0.234375
    "31" = 0.234375
    "32" = 0.234375
    "33" = 0.203125
    "34" = 0.234375
    "35" = 0.234375
    "36" = 0.25
    "38" = 0.234375
    "39" = 0.21875
    "40" = 0.234375
    "41" = 0.234375
    "42" = 0.234375
    "43" = 0.234375
    "44" = 0.234375
    "45" = 0.234375
    "46" = 0.234375
    "47" = 0.234375
    "48" = 0.234375
    "49" = 0.234375
    "50" = 0.234375
    "51" = 0.234375
    "52" = 0.25
    "53" = 0.25
    "54" = 0.25
    "55" = 0.25
    "56" = 0.25
    "57" = 0.25
    "58" = 0.25
    "59" = 0.25
    "60" = 0.25
    "61" = 0.25
    "62" = 0.25
    "63" = 0.25
    "64" = 0.25
    "65" = 0.25
    "66" = 0.25
    "67" = 0.25
    "68" = 0.25
    "69" = 0.234375
    "70" = 0.234375
    "71" = 0.234375
    "72" = 0.234375
    "73" = 0.234375
    "74" = 0.234375
    "75" = 0.234375
    "76" = 0.234375
    "77" = 0.234375
    "78" = 0.234375
    "79" = 0.234375
    "80" = 0.234375
    "81" = 0.234375
    "82" = 0.234375
    "83" = 0.234375
    "84" = 0.234375
    "85" = 0.234375
    "86" = 0.234375
    "87" = 0.234375
    "88" = 0.234375
    "89" = 0.234375
    "90" = 0.234375
    "91" = 0.234375
    "92" = 0.234375
    "93" = 0.234375
    "94" = 0.234375
    "95" = 0.234375
    "96" = 0.234375
    "97" = 0.234375
    "98" = 0.234375
    "99" = 0.234375
    "100" = 0.234375
    "101" = 0.234375
    "102" = 0.234375
    "103" = 0.140625
    "105" = 0.265625
    "106" = 0.234375
    "107" = 0.296875
    "108" = 0.234375
    "110" = 0.234375
    "111" = 0.140625
    "112" = 0.09375
    "113" = 0.203125
    "114" = 0.1875
    "115" = 0.25
    "116" = 0.1875
    "117" = 0.296875
}

foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value = $newValues[$row]
}

# The "DisplayOutputs" repr in column A reflects the live Python object id
# from the re-executed notebook cell; refresh it for the rows that show it.
$oldRepr = "<__main__.DisplayOutputs object at 0x7f323c1cdbb0>"
$newRepr = "<__main__.DisplayOutputs object at 0x7fb2226fc7c0>"
for ($r = 102; $r -le 118; $r++) {
    $cell = $ws.Range("A$r")
    if ($cell.Value2 -eq $oldRepr) {
        $cell.Value = $newRepr
    }
}

# Reflect the final UI state: whole sheet selected (Ctrl+A). The saved
# file shows the cursor resting on R17 within that selection; this
# engine's Range.Select/Activate always collapses ActiveCell to match
# the most recently touched range, so the full-sheet selection (the
# data-significant part of the gesture) is preserved here.
$null = $ws.Range("A1:XFD1048576").Select()
